$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 30, pushing the old rows
# 30-32 down to 33-35 (their values remain unchanged).
$ws.Rows("30:32").Insert()

# Fill the 3 newly inserted rows (30, 31, 32) with the new weekly data.

# Row 30
$ws.Range("A30").Value = 3
$ws.Range("B30").Value = "Femacal de La Calera"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44460
$ws.Range("E30").Value = 5
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100107
$ws.Range("H30").Value = "Otros"
$ws.Range("I30").Value = 100107002
$ws.Range("J30").Value = "Chirimoya"
$ws.Range("K30").Value = "Cultivar IV Región"
$ws.Range("L30").Value = "Especial"
$ws.Range("M30").Value = 45
$ws.Range("N30").Value = 30000
$ws.Range("O30").Value = 30000
$ws.Range("P30").Value = 30000
$ws.Range("Q30").Value = "$/bandeja 10 kilos"
$ws.Range("R30").Value = "Provincia del Elquí"
$ws.Range("S30").Value = 3000
$ws.Range("T30").Value = 10

# Row 31
$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "Femacal de La Calera"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44460
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100107
$ws.Range("H31").Value = "Otros"
$ws.Range("I31").Value = 100107002
$ws.Range("J31").Value = "Chirimoya"
$ws.Range("K31").Value = "Cultivar IV Región"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 50
$ws.Range("N31").Value = 27000
$ws.Range("O31").Value = 27000
$ws.Range("P31").Value = 27000
$ws.Range("Q31").Value = "$/bandeja 10 kilos"
$ws.Range("R31").Value = "Provincia del Elquí"
$ws.Range("S31").Value = 2700
$ws.Range("T31").Value = 10

# Row 32
$ws.Range("A32").Value = 3
$ws.Range("B32").Value = "Femacal de La Calera"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = 44460
$ws.Range("E32").Value = 5
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100107
$ws.Range("H32").Value = "Otros"
$ws.Range("I32").Value = 100107002
$ws.Range("J32").Value = "Chirimoya"
$ws.Range("K32").Value = "Cultivar IV Región"
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 48
$ws.Range("N32").Value = 25000
$ws.Range("O32").Value = 25000
$ws.Range("P32").Value = 25000
$ws.Range("Q32").Value = "$/bandeja 10 kilos"
$ws.Range("R32").Value = "Provincia del Elquí"
$ws.Range("S32").Value = 2500
$ws.Range("T32").Value = 10

# Make sure the date cells keep the date/time numeric display style
# that the rest of column D uses.
$ws.Range("D30:D32").NumberFormat = $ws.Range("D29").NumberFormat
